# Insert two new rows of weekly price data at the top of the Betarraga data
# block (row 290), pushing the existing rows down by two. The workbook's
# dimension grows from A1:R388 to A1:R390.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new blank rows before row 290 - this shifts rows 290:388 down to 292:390
$ws.Range("A290:R291").EntireRow.Insert()

# Give the date cells in the two new rows the same number format used by the
# rest of column D (style index 2 / numFmtId 165, "YYYY-MM-DD HH:MM:SS").
$ws.Range("D290:D291").NumberFormat = $ws.Range("D292").NumberFormat

# Row 290 - "Primera" quality entry for the newest week (2021-09-29)
$ws.Cells.Item(290, 1).Value = 3
$ws.Cells.Item(290, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(290, 3).Value = "Coquimbo"
$ws.Cells.Item(290, 4).Value = 44468
$ws.Cells.Item(290, 5).Value = 5
$ws.Cells.Item(290, 6).Value = 100114014
$ws.Cells.Item(290, 7).Value = "Betarraga"
$ws.Cells.Item(290, 8).Value = "Sin especificar"
$ws.Cells.Item(290, 9).Value = "Primera"
$ws.Cells.Item(290, 10).Value = 3100
$ws.Cells.Item(290, 11).Value = 500
$ws.Cells.Item(290, 12).Value = 550
$ws.Cells.Item(290, 13).Value = 524
$ws.Cells.Item(290, 14).Value = "`$/paquete 4 unidades"
$ws.Cells.Item(290, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(290, 16).Value = 131
$ws.Cells.Item(290, 17).Value = 4
$ws.Cells.Item(290, 18).Value = "Hortaliza"

# Row 291 - "Segunda" quality entry for the newest week (2021-09-29)
$ws.Cells.Item(291, 1).Value = 3
$ws.Cells.Item(291, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(291, 3).Value = "Coquimbo"
$ws.Cells.Item(291, 4).Value = 44468
$ws.Cells.Item(291, 5).Value = 5
$ws.Cells.Item(291, 6).Value = 100114014
$ws.Cells.Item(291, 7).Value = "Betarraga"
$ws.Cells.Item(291, 8).Value = "Sin especificar"
$ws.Cells.Item(291, 9).Value = "Segunda"
$ws.Cells.Item(291, 10).Value = 1400
$ws.Cells.Item(291, 11).Value = 400
$ws.Cells.Item(291, 12).Value = 400
$ws.Cells.Item(291, 13).Value = 400
$ws.Cells.Item(291, 14).Value = "`$/paquete 4 unidades"
$ws.Cells.Item(291, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(291, 16).Value = 100
$ws.Cells.Item(291, 17).Value = 4
$ws.Cells.Item(291, 18).Value = "Hortaliza"
